$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") values for rows 2-28 were all 45422 (2024-05-10)
# and are being updated to 45423 (2024-05-11).
for ($r = 2; $r -le 28; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45422) {
        $cell.Value2 = 45423
    }
}
